# Applies the "recategorisation des citations" edit to metadata_avis.xlsx
# - fills in the "citations_autres_avis" (column Q) for a batch of rows
# - fixes two "divergence" (column N) booleans from False to True
# - un-hides / resizes column C
# - updates the saved selection / scroll position on the metadata sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metadata")

# --- Column Q ("citations_autres_avis") fill-in, entered bottom-to-top ---
$ws.Range("Q101").Value = 39
$ws.Range("Q100").Value = '4;5;17;25;30'
$ws.Range("Q96").Value = 49
$ws.Range("Q95").Value = '39;40'
$ws.Range("Q94").Value = 1
$ws.Range("Q93").Value = '1;23'
$ws.Range("Q92").Value = '8;9;16;21'
$ws.Range("Q89").Value = '46;42;24'
$ws.Range("Q88").Value = '34;39;43;38'
$ws.Range("Q87").Value = '57;58'
$ws.Range("Q86").Value = '1;52;53;54;19;8;56;46;4;9;25;38;35;36;41;39;40;43;44;45;49;50;51;54'
$ws.Range("Q84").Value = 57
$ws.Range("Q83").Value = '26;42;58;57;59'
$ws.Range("Q81").Value = '24;42;60'
$ws.Range("Q79").Value = '60;54;53'
$ws.Range("Q78").Value = '5;25;37;46;19;42;65;57'
$ws.Range("Q77").Value = 56
$ws.Range("Q76").Value = '58;46'
$ws.Range("Q75").Value = '58;70'
$ws.Range("Q74").Value = '8;19;60;70'
$ws.Range("Q73").Value = 2
$ws.Range("Q71").Value = '72;42'
$ws.Range("Q70").Value = '25;70;72'
$ws.Range("Q68").Value = '41;57;17'
$ws.Range("Q67").Value = '58;73;55;70;71;2;12;57;38'
$ws.Range("Q66").Value = 46
$ws.Range("Q65").Value = '45;35;75'
$ws.Range("Q63").Value = '76;46'
$ws.Range("Q61").Value = 79
$ws.Range("Q59").Value = '86;63;65;70;84'
$ws.Range("Q46").Value = 90
$ws.Range("Q17").Value = '112;116;109;77;127;105;46;124;120;115;122;126;90;113;26;63;87;108;121;128'

# --- "degree in or out" (divergence) fixes ---
$ws.Range("N58").Value = $true
$ws.Range("N82").Value = $true

# --- column C: un-hide and resize ---
$ws.Columns.Item(3).Hidden = $false
$ws.Columns.Item(3).ColumnWidth = 9.5703125

# --- restore view / selection state ---
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("Q17").Select()
